$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update task text for the "Image-viewer-actions" block (rows 19-21) ---
$ws.Range("B19").Value = "Click on the image to add a comment"
$ws.Range("B20").Value = "Select the geometic option comment"

$pad = "".PadLeft(1194)
$liveTrace = "Select the live trace option on the" + $pad + "comment"
$ws.Range("B21").Value = $liveTrace

# --- Insert a new 4th task row (row 22) for the pin-comment task ---
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null
$ws.Range("B21").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null

$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Select the pin option, add some text and create a comment"
$ws.Rows(22).RowHeight = 34

# --- Fix "dificult" typo -> "difficult" in the three scale questions ---
$ws.Range("A24").Value = 'In scale 1-7, being 1 very difficult and 7 very easy how easy was to perform the "General-flow-test" tasks?'
$ws.Range("A25").Value = 'In scale 1-7, being 1 very difficult and 7 very easy how easy was to perform the "Asset-explorer-actions" tasks?'
$ws.Range("A26").Value = 'In scale 1-7, being 1 very difficult and 7 very easy how easy was to perform the "Image-viewer-actions" tasks?'

# --- Update the view: drop the stale topLeftCell/selection and focus B14 ---
$ws.Range("B14").Select() | Out-Null
